$d = $word.ActiveDocument

# RGB(201, 33, 30) == hex C9211E, encoded as an OLE color (0x00BBGGRR,
# i.e. R + G*256 + B*65536).
$calloutColor = 201 + (33 * 256) + (30 * 65536)

# Locate the paragraph that contains the target sentence.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*This might cause some frustration for the users.*") {

        # Add a trailing space after "users." so the sentence reads
        # "...frustration for the users. " (xml:space="preserve").
        [void]$p.Range.Find.Execute("users.", $false, $false, $false, $false, $false, `
                               $true, 1, $false, "users. ", 2)

        # Append "– " in red (C9211E), mirroring the existing callout style
        # already used elsewhere in this document.
        $dash = $p.Range
        $dash.SetRange($dash.End - 1, $dash.End - 1)
        $dash.InsertAfter([string][char]0x2013 + " ")
        $dash.Font.Color = $calloutColor
        $dash.LanguageID = "en-US"

        # Append "Done (Observation header and text retain)" in the same red.
        $done = $p.Range
        $done.SetRange($done.End - 1, $done.End - 1)
        $done.InsertAfter("Done (Observation header and text retain)")
        $done.Font.Color = $calloutColor
        $done.LanguageID = "en-US"

        break
    }
}
